$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("setpoints")

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1

$ws.Range("E3").Select()
